$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MSI")

# Row 4: Inventory
$ws.Range("B4").Value = 508000000.0
$ws.Range("C4").Value = 489000000.0
$ws.Range("D4").Value = 449000000.0
$ws.Range("E4").Value = 442000000.0
$ws.Range("F4").Value = 447000000.0

# Row 15: Accounts Payable
$ws.Range("B15").Value = 612000000.0
$ws.Range("C15").Value = 536000000.0
$ws.Range("D15").Value = 498000000.0
$ws.Range("E15").Value = 531000000.0
$ws.Range("F15").Value = 618000000.0

# Row 22: Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = -786000000.0
$ws.Range("C22").Value = -697000000.0
$ws.Range("D22").Value = -747000000.0
$ws.Range("E22").Value = -747000000.0
$ws.Range("F22").Value = -759000000.0
